# Fix wrong dimension name in documentation (miterBoxHeight -> miterBoxWidth)
# and bump the auto-date placeholder text on every layout + the master
# from 2021-02-19 to 2021-02-21 (file was re-saved two days later).

$p = $ppt.ActivePresentation

function Set-FullRangeText($textRange, [string]$newText) {
    # Replacing the whole range starting at character 1 (rather than
    # assigning TextRange.Text directly) keeps PowerPoint from re-diffing
    # run boundaries, so the result collapses cleanly onto the first run's
    # formatting instead of being split into multiple runs.
    $len = $textRange.Length
    if ($len -gt 0) {
        $full = $textRange.Characters(1, $len)
        $full.Text = $newText
    } else {
        $textRange.Text = $newText
    }
}

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2021-02-19") {
                Set-FullRangeText $tr "2021-02-21"
            }
        }
    }
}

# --- Slide master date placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout's date placeholder ---
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# --- Fix the "miterBoxHeight" typo to "miterBoxWidth" on slide 8 ---
$s8 = $p.Slides.Item(8)
for ($i = 1; $i -le $s8.Shapes.Count; $i++) {
    $sh = $s8.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "miterBoxHeight") {
            Set-FullRangeText $tr "miterBoxWidth"
        }
    }
}
